$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6577333333333334
$ws.Range("H2").Value = 1.9732
$ws.Range("M2").Value = 1.665504333333333
$ws.Range("N2").Value = 4.996513
$ws.Range("Q2").Value = 1.095457716844445
$ws.Range("R2").Value = 9.859119451600002
